# AT.docx — generalize "CivicActions Information Security Office" / "ISSO" /
# "CivicActions Security" wording to "CivicActions' Security Office", plus a
# few accompanying wording tweaks, per the commit's new 'Contractor'
# component generalization.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

# "Both regular and ad hoc training..." paragraph (3 runs touched).
Replace-Text `
    "must complete Security Awareness trainings upon being hired and at least annually thereafter. CivicActions Operations will not create accounts" `
    "must complete Security Awareness training upon being hired and at least annually thereafter. CivicActions Operations staff will not create accounts"

Replace-Text `
    "Results from this survey are compiled by the Director of Human Resources and used to refine future training efforts." `
    "Results from this survey are compiled by the Office of Human Resources and used to refine future training efforts."

Replace-Text `
    "Ad Hoc Security Awareness: The CivicActions ISSO oversees the approximately bi-monthly distribution of security awareness tips and articles to the all CivicActions employees." `
    "Ad Hoc Security Awareness: The CivicActions’ Security Office oversees the approximately bi-monthly distribution of security awareness tips and articles to all CivicActions employees."

# "In the event of a major system change..." paragraph.
Replace-Text `
    "Specific training type, medium and delivery method is dependent upon the nature of the system change." `
    "Specific training types, mediums, and delivery methods are dependent upon the nature of the system change."

# "CivicActions personnel with security responsibilities..." paragraph.
Replace-Text `
    "The CivicActions ISSO is responsible for creating the content of the training. The role-based training is provided and tracked by the CivicActions Information Security Office." `
    "The CivicActions’ Security Office is responsible for creating the content of the training. The role-based training is provided and tracked by the CivicActions’ Security Office."

# "The Project manager in collaboration with CivicActions Security determines..." paragraph.
Replace-Text `
    "The Project manager in collaboration with CivicActions Security determines whether a change to the information system requires any modifications and updates to the security awareness training program and if so, works with the CivicActions Security to implement the change." `
    "The Project Manager in collaboration with CivicActions’ Security Office determines whether a change to the information system requires any modifications and updates to the security awareness training program and if so, works with the CivicActions’ Security Office to implement the change."

# "CivicActions’ Security provides users with security responsibilities..." paragraph.
Replace-Text `
    "CivicActions’ Security provides users with security responsibilities role-based security training on an annual basis. The training is provided and tracked by the CivicActions Information Security Office." `
    "CivicActions’ Security Office provides users with security responsibilities role-based security training on an annual basis. The training is provided and tracked by the CivicActions’ Security Office."

# "The CivicActions Information Security Office tracks all security awareness training..." paragraph.
Replace-Text `
    "The CivicActions Information Security Office tracks all security awareness training within the organization and ensures that all employees have successfully completed training when required. The training records are stored and tracked in a spreadsheet maintained by the CivicActions Information Security Office." `
    "The CivicActions’ Security Office tracks all security awareness training within the organization and ensures that all employees have successfully completed training when required. The training records are stored and tracked in a spreadsheet maintained by the CivicActions’ Security Office."

# "Training records are tracked and maintained by..." paragraph.
Replace-Text `
    "Training records are tracked and maintained by the CivicActions Information Security Office. Records are maintained permanently." `
    "Training records are tracked and maintained by the CivicActions’ Security Office. Records are maintained permanently."
